$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column K ("LastRound"/"特殊效果" columns all shift
# one column to the right: K->L, L->M) to make room for the new "Price" field.
$ws.Columns("K").Insert()

# Give the freshly inserted column K the same look as the neighbouring
# "属性加成比例" columns (G:J): centred horizontally + vertically, matching
# cell style used throughout the rest of the header block.
$ws.Range("K1:K2").HorizontalAlignment = -4108
$ws.Range("K1:K2").VerticalAlignment = -4108

# Row 1/2 headers (merged-comment style block) + row 3 field-name row.
$ws.Range("K1").Value = "购买价格"
$ws.Range("K2").Value = "int"
$ws.Range("K3").Value = "Price"

# Data rows: purchase price for each of the three items.
$ws.Range("K4").Value = 20
$ws.Range("K5").Value = 30
$ws.Range("K6").Value = 50

# Match the column width used by the G:J group (14.5398230088496 in raw
# OOXML character-width units == ColumnWidth 13.825537294563887 once the
# standard 5px/MDW padding is backed out).
$ws.Range("K1").ColumnWidth = 13.825537294563887

# Highlight (in red) the leading "第一、二行" / "第三行此列之前" runs of the
# explanatory comment cells (now in column M) so they stand out like the
# other emphasised fragments already in those same rich-text strings.
$ws.Range("M1").Characters(1, 5).Font.Color = 255
$ws.Range("M2").Characters(1, 7).Font.Color = 255

# Leave the cursor where the author's edit session ended up.
[void]$ws.Range("K6").Select()
